$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target values ("12", "21", ...) look numeric, so a plain .Value
# assignment on a General-formatted cell would be auto-coerced to a
# number by Excel. The source data stores these as text, so force text
# entry (temporarily mark the cell as Text) and then restore the
# cell's style to Normal so no stray formatting is left behind -
# only the cell *value* should change, matching the original diff.
$rng = $ws.Range("A2:H2")
$rng.NumberFormat = "@"

$ws.Range("A2").Value = "12"
$ws.Range("B2").Value = "12"
$ws.Range("C2").Value = "12"
$ws.Range("D2").Value = "12"
$ws.Range("E2").Value = "12"
$ws.Range("F2").Value = "12"
$ws.Range("G2").Value = "21"
$ws.Range("H2").Value = "12"

$rng.Style = "Normal"
